$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the numeric-looking "Price" cells as Text so assigning their
# updated values keeps the exact original string representation (e.g.
# "1.000", "0.07816") instead of Excel auto-coercing them into numbers.
$textCells = @("D4", "D5", "D6", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D17", "D19", "D21", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D41", "D42", "D43", "D44", "D45", "D49", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply the updated crypto market values
$ws.Range("D2").Value = "29.247.76"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "1.863.50"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "0.7057"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").Value = "242.29"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "0.07816"
$ws.Range("E8").Value = "  -2.10%  "
$ws.Range("D9").Value = "0.3114"
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").Value = "24.30"
$ws.Range("E10").Value = "  -3.70%  "
$ws.Range("D11").Value = "0.08004"
$ws.Range("E11").Value = "  -4.01%  "
$ws.Range("D12").Value = "1.861.49"
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("D13").Value = "93.66"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("D14").Value = "5.175"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").Value = "0.6964"
$ws.Range("E15").Value = "  -3.31%  "
$ws.Range("D16").Value = "6.356"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.000008307"
$ws.Range("E17").Value = "  -2.13%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "29.195.79"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("D19").Value = "252.75"
$ws.Range("E19").Value = "  +4.91%  "
$ws.Range("D20").Value = "2.142.21"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").Value = "13.10"
$ws.Range("E21").Value = "  -1.20%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "7.503"
$ws.Range("E23").Value = "  -4.28%  "
$ws.Range("D25").Value = "0.1560"
$ws.Range("E25").Value = "  -1.76%  "
$ws.Range("D26").Value = "8.996"
$ws.Range("E26").Value = "  -0.84%  "
$ws.Range("D27").Value = "159.46"
$ws.Range("E27").Value = "  -2.64%  "
$ws.Range("D28").Value = "18.81"
$ws.Range("E28").Value = "  +1.31%  "
$ws.Range("D29").Value = "1.500"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("D30").Value = "4.288"
$ws.Range("E30").Value = "  -2.89%  "
$ws.Range("D31").Value = "4.269"
$ws.Range("E31").Value = "  -1.60%  "
$ws.Range("D32").Value = "1.210"
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("D33").Value = "0.05260"
$ws.Range("E33").Value = "  -2.06%  "
$ws.Range("D34").Value = "1.890"
$ws.Range("E34").Value = "  -3.11%  "
$ws.Range("D35").Value = "0.7450"
$ws.Range("E35").Value = "  -0.49%  "
$ws.Range("D36").Value = "1.158"
$ws.Range("E36").Value = "  -2.09%  "
$ws.Range("D37").Value = "2.706"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "0.01860"
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("D39").Value = "1.257.33"
$ws.Range("E39").Value = "  -2.41%  "
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("D41").Value = "6.298"
$ws.Range("E41").Value = "  -4.59%  "
$ws.Range("D42").Value = "0.8999"
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("D43").Value = "111.10"
$ws.Range("E43").Value = "  -0.45%  "
$ws.Range("D44").Value = "71.75"
$ws.Range("E44").Value = "  -2.43%  "
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").Value = "2.036.64"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D49").Value = "1.781"
$ws.Range("E49").Value = "  -1.55%  "
$ws.Range("D50").Value = "9.391"
$ws.Range("E50").Value = "  -1.06%  "
$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").Value = "0.4302"
$ws.Range("E51").Value = "  -2.10%  "
